# Swap columns C ("codeforiati:group-name") and D ("codeforiati:group-code")
# so that the group-code column comes before the group-name column,
# matching the reordering performed upstream in codeforIATI/codelists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value()
    $dValue = $dCell.Value()

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
